$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a new column at R (column 18), shifting the old R column (and
# everything after it) one column to the right so it becomes column S.
$ws.Columns.Item(18).Insert()

# The header cell that used to live in Q1 ("pricing_interest_rate_type")
# is duplicated into the freshly inserted R1 header cell.
$ws.Range("R1").Value = $ws.Range("Q1").Value()

# Fill in the new R column (rows 2-11) with the values that previously
# lived in Q (the row index / numeric counters), and update a handful of
# Q cells with the new values from the merge.
$ws.Range("R2").Value = 1
$ws.Range("R3").Value = 2
$ws.Range("R4").Value = 3
$ws.Range("R5").Value = 4
$ws.Range("R6").Value = 5
$ws.Range("R7").Value = 6
$ws.Range("R8").Value = 7
$ws.Range("R9").Value = 900
$ws.Range("R10").Value = 1
$ws.Range("R11").Value = 3

$ws.Range("Q8").Value = 999
$ws.Range("Q9").Value = 0
$ws.Range("Q10").Value = 10
$ws.Range("Q11").Value = 1000

# The Q column no longer carries the wrap-text style it used to (that
# style now belongs to the newly inserted R column).
$ws.Range("Q1:Q11").Style = "Normal"

# Update the remembered selection saved with the sheet.
$ws.Range("Q16").Select()
